# "04_02 Mapping Study" update
#
# The author re-opened the workbook (on a different machine/Excel build)
# and updated the paper counts on the "Ausgangstabelle RQ1" sheet, then
# left the selection/zoom in a different spot and touched the page setup
# for that sheet before saving again.
#
# (The bulk of the underlying XML diff - fileVersion/rupBuild, absPath,
# revisionPtr GUIDs, bookViews window geometry, sub-pixel column-width
# jitter, and auto row-height jitter tied to a Mac/Windows font-metric
# difference - is incidental re-save noise from the different Excel build
# and carries no editable content; it is not reproduced here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ausgangstabelle RQ1")

# --- updated paper counts -------------------------------------------------
$ws.Range("B2").Value = 9
$ws.Range("C2").Value = 7

$ws.Range("B3").Value = 15
$ws.Range("C3").Value = 13
$ws.Range("D3").Value = 12
$ws.Range("E3").Value = 12

$ws.Range("B4").Value = 13
$ws.Range("C4").Value = 11
$ws.Range("D4").Value = 14
$ws.Range("E4").Value = 10

$ws.Range("B5").Value = 15
$ws.Range("C5").Value = 13
$ws.Range("E5").Value = 12
$ws.Range("F5").Value = 11

# --- view state: zoom in and move the selection ---------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 125
$ws.Range("E15").Select()

# --- page setup for this sheet ---------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
